$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.631.61'
$ws.Range("E2").Value = '  +0.93%  '

$ws.Range("D3").Value = '1.563.63'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("D5").Value = '''210.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '

$ws.Range("D6").Value = '''0.521'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.33%  '

$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("D8").Value = '''24.76'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.12%  '

$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("D10").Value = '''0.0588'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.15%  '

$ws.Range("D11").Value = '''0.0900'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("D12").Value = '1.787.37'
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").Value = '1.567.47'
$ws.Range("E13").Value = '  -0.02%  '

$ws.Range("D14").Value = '28.667.62'
$ws.Range("E14").Value = '  +1.12%  '

$ws.Range("D15").Value = '''0.516'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("E16").Value = '  -0.79%  '

$ws.Range("D17").Value = '''61.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.28%  '

$ws.Range("D18").Value = '''227.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.38%  '

$ws.Range("D19").Value = '''7.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.18%  '

$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("E22").Value = '  -0.67%  '

$ws.Range("D23").Value = '''9.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.26%  '

$ws.Range("D24").Value = '''2.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").Value = '''151.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("E26").Value = '  +3.70%  '

$ws.Range("E27").Value = '  -0.65%  '

$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '

$ws.Range("E29").Value = '  -1.00%  '

$ws.Range("E30").Value = '  -3.78%  '

$ws.Range("E31").Value = '  -0.89%  '

$ws.Range("D32").Value = '''3.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.16%  '

$ws.Range("D33").Value = '1.402.69'
$ws.Range("E33").Value = '  +1.27%  '

$ws.Range("E34").Value = '  -2.40%  '

$ws.Range("E35").Value = '  -3.71%  '

$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("E37").Value = '  +2.99%  '

$ws.Range("E38").Value = '  -2.08%  '

$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("D40").Value = '''0.518'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("E41").Value = '  -0.88%  '

$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("E43").Value = '  -2.28%  '

$ws.Range("D44").Value = '''0.0460'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.67%  '

$ws.Range("D45").Value = '''63.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.72%  '

$ws.Range("E46").Value = '  -2.36%  '

$ws.Range("D47").Value = '1.699.36'
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").Value = '''0.840'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.73%  '

$ws.Range("D49").Value = '''84.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.84%  '

$ws.Range("D50").Value = '''42.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.76%  '

$ws.Range("D51").Value = '''0.0511'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '
